$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 95441.39999999999
$ws.Range("J3").Value = 95441.39999999999
$ws.Range("L3").Value = 95441.39999999999
$ws.Range("N3").Value = -95669.39999999999
$ws.Range("H17").Value = 428.06
$ws.Range("J17").Value = 428.06
$ws.Range("L17").Value = 1284.18
$ws.Range("N17").Value = -1620.18
$ws.Range("H19").Value = 2891.6667
$ws.Range("I19").Value = 2349.5
$ws.Range("K19").Value = 2349.5
$ws.Range("M19").Value = -2174.5
$ws.Range("H32").Value = 2499.3333
$ws.Range("J32").Value = 2999.5
$ws.Range("L32").Value = 2999.5
$ws.Range("N32").Value = -3651.5
$ws.Range("H86").Value = 14577.077
$ws.Range("I86").Value = 14050.3
$ws.Range("J86").Value = 16333
$ws.Range("K86").Value = 14050.3
$ws.Range("L86").Value = 16333
$ws.Range("M86").Value = -12927.3
$ws.Range("N86").Value = -18579
$ws.Range("H89").Value = 14577.077
$ws.Range("I89").Value = 14050.3
$ws.Range("J89").Value = 16333
$ws.Range("K89").Value = 70251.5
$ws.Range("L89").Value = 81665
$ws.Range("M89").Value = -64635.5
$ws.Range("N89").Value = -92897
$ws.Range("H102").Value = 95441.39999999999
$ws.Range("J102").Value = 95441.39999999999
$ws.Range("L102").Value = 95441.39999999999
$ws.Range("N102").Value = -101931.4
$ws.Range("H116").Value = 7050.1875
$ws.Range("I116").Value = 6483.0835
$ws.Range("K116").Value = 6483.0835
$ws.Range("M116").Value = -3041.0835
$ws.Range("H132").Value = 2522.2273
$ws.Range("I132").Value = 2109.5386
$ws.Range("K132").Value = 6328.6158
$ws.Range("M132").Value = -3798.6158
$ws.Range("H138").Value = 2202.9285
$ws.Range("I138").Value = 1295.3077
$ws.Range("J138").Value = 2609.7932
$ws.Range("K138").Value = 3885.9231
$ws.Range("L138").Value = 7829.3796
$ws.Range("M138").Value = 1254.0769
$ws.Range("N138").Value = -18109.3796
$ws.Range("H141").Value = 2599
$ws.Range("I141").Value = 2599
$ws.Range("K141").Value = 7797
$ws.Range("M141").Value = -2617

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 12751.81
$ws.Range("I102").Value = 21884.143
$ws.Range("J102").Value = 8185.643
$ws.Range("K102").Value = 21884.143
$ws.Range("L102").Value = 8185.643
$ws.Range("M102").Value = -20262.143
$ws.Range("N102").Value = -11429.643
$ws.Range("H110").Value = 1647.1482
$ws.Range("I110").Value = 1632.25
$ws.Range("J110").Value = 1766.3334
$ws.Range("K110").Value = 1632.25
$ws.Range("L110").Value = 1766.3334
$ws.Range("M110").Value = 412.75
$ws.Range("N110").Value = -5856.3334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5108.364
$ws.Range("I20").Value = 4973.75
$ws.Range("J20").Value = 5467.3335
$ws.Range("K20").Value = 4973.75
$ws.Range("L20").Value = 5467.3335
$ws.Range("M20").Value = -4726.75
$ws.Range("N20").Value = -5961.3335
$ws.Range("H99").Value = 2556.5715
$ws.Range("I99").Value = 1750.5
$ws.Range("K99").Value = 1750.5
$ws.Range("M99").Value = -252.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2394.5557
$ws.Range("I16").Value = 2401.6
$ws.Range("J16").Value = 2385.75
$ws.Range("K16").Value = 2401.6
$ws.Range("L16").Value = 2385.75
$ws.Range("M16").Value = -2114.6
$ws.Range("N16").Value = -2959.75
$ws.Range("H31").Value = 461682.1
$ws.Range("I31").Value = 3154.0386
$ws.Range("K31").Value = 3154.0386
$ws.Range("M31").Value = -2859.0386
$ws.Range("H34").Value = 461682.1
$ws.Range("I34").Value = 3154.0386
$ws.Range("K34").Value = 3154.0386
$ws.Range("M34").Value = -2952.0386
$ws.Range("H105").Value = 2295.9285
$ws.Range("I105").Value = 1790.7142
$ws.Range("J105").Value = 2801.1428
$ws.Range("K105").Value = 1790.7142
$ws.Range("L105").Value = 2801.1428
$ws.Range("M105").Value = -43.71419999999989
$ws.Range("N105").Value = -6295.1428
$ws.Range("H113").Value = 2394.5557
$ws.Range("I113").Value = 2401.6
$ws.Range("J113").Value = 2385.75
$ws.Range("K113").Value = 2401.6
$ws.Range("L113").Value = 2385.75
$ws.Range("M113").Value = -231.5999999999999
$ws.Range("N113").Value = -6725.75
$ws.Range("H132").Value = 4306.3
$ws.Range("I132").Value = 4256.125
$ws.Range("K132").Value = 12768.375
$ws.Range("M132").Value = -10238.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 8000
$ws.Range("I42").Value = 8000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 24000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -23466
$ws.Range("N42").ClearContents()
$ws.Range("H131").Value = 7707.2666
$ws.Range("J131").Value = 6121
$ws.Range("L131").Value = 18363
$ws.Range("N131").Value = -28443
$ws.Range("H134").Value = 4130.8184
$ws.Range("I134").Value = 2680.0625
$ws.Range("K134").Value = 8040.1875
$ws.Range("M134").Value = -2970.1875
$ws.Range("H140").Value = 217379.58
$ws.Range("I140").Value = 275560.62
$ws.Range("K140").Value = 826681.86
$ws.Range("M140").Value = -821501.86

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8019.6665
$ws.Range("I102").Value = 4623.6
$ws.Range("K102").Value = 4623.6
$ws.Range("M102").Value = -3001.6
$ws.Range("H122").Value = 2036.125
$ws.Range("I122").Value = 1715
$ws.Range("K122").Value = 5145
$ws.Range("M122").Value = -2695

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3137.1052
$ws.Range("I22").Value = 3266.9167
$ws.Range("J22").Value = 2914.5715
$ws.Range("K22").Value = 3266.9167
$ws.Range("L22").Value = 2914.5715
$ws.Range("M22").Value = -2971.9167
$ws.Range("N22").Value = -3504.5715
$ws.Range("H27").Value = 3137.1052
$ws.Range("I27").Value = 3266.9167
$ws.Range("J27").Value = 2914.5715
$ws.Range("K27").Value = 3266.9167
$ws.Range("L27").Value = 2914.5715
$ws.Range("M27").Value = -3159.9167
$ws.Range("N27").Value = -3128.5715
$ws.Range("H46").Value = 4458.381
$ws.Range("J46").Value = 4665
$ws.Range("L46").Value = 4665
$ws.Range("N46").Value = -5041
$ws.Range("I55").Value = 100000264
$ws.Range("K55").Value = 100000264
$ws.Range("M55").Value = -100000091
$ws.Range("H68").Value = 1077.6
$ws.Range("I68").Value = 1077.6
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1077.6
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -328.5999999999999
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 1077.6
$ws.Range("I71").Value = 1077.6
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 5388
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -1644
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 51914.668
$ws.Range("J74").Value = 67872
$ws.Range("L74").Value = 67872
$ws.Range("N74").Value = -69868
$ws.Range("H77").Value = 51914.668
$ws.Range("J77").Value = 67872
$ws.Range("L77").Value = 203616
$ws.Range("N77").Value = -213600
$ws.Range("H132").Value = 145341.4
$ws.Range("I132").Value = 112901.78
$ws.Range("J132").Value = 194000.83
$ws.Range("K132").Value = 338705.34
$ws.Range("L132").Value = 582002.49
$ws.Range("M132").Value = -336175.34
$ws.Range("N132").Value = -587062.49

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4786.057
$ws.Range("I122").Value = 2294.25
$ws.Range("K122").Value = 6882.75
$ws.Range("M122").Value = -4432.75
